$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = 111936868
$ws.Range("B11").Value = 89423
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = "Granticka"
$ws.Range("G11").Value = "Porodaedalea chrysoloma"
$ws.Range("H11").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("K11").ClearContents()
$ws.Range("L11").ClearContents()
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()
$ws.Range("Q11").Value = 448988.017639213
$ws.Range("R11").Value = 7087186.778340456
$ws.Range("AC11").ClearContents()

# Row 12
$ws.Range("A12").Value = 111936870
$ws.Range("B12").Value = 89423
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = "Granticka"
$ws.Range("G12").Value = "Porodaedalea chrysoloma"
$ws.Range("H12").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q12").Value = 449019.027096529
$ws.Range("R12").Value = 7087276.979166135

# Row 13
$ws.Range("A13").Value = 111936865
$ws.Range("B13").Value = 89423
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 5432
$ws.Range("F13").Value = "Granticka"
$ws.Range("G13").Value = "Porodaedalea chrysoloma"
$ws.Range("H13").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q13").Value = 448738.4239939091
$ws.Range("R13").Value = 7087426.42220111

# Row 14
$ws.Range("A14").Value = 111936893
$ws.Range("B14").Value = 77515
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = "Garnlav"
$ws.Range("G14").Value = "Alectoria sarmentosa"
$ws.Range("H14").Value = "(Ach.) Ach."
$ws.Range("K14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()
$ws.Range("Q14").Value = 448742.3011697636
$ws.Range("R14").Value = 7087501.648173723
$ws.Range("AC14").ClearContents()

# Row 15
$ws.Range("A15").Value = 111936866
$ws.Range("Q15").Value = 448765.5992023234
$ws.Range("R15").Value = 7087416.731054713

# Row 16
$ws.Range("A16").Value = 111936798
$ws.Range("B16").Value = 56398
$ws.Range("E16").Value = 100109
$ws.Range("F16").Value = "Tretåig hackspett"
$ws.Range("G16").Value = "Picoides tridactylus"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("Q16").Value = 448923.1251473558
$ws.Range("R16").Value = 7087371.00725084
$ws.Range("AC16").Value = "ringhack äldre"

# Row 17
$ws.Range("A17").Value = 111936792
$ws.Range("B17").Value = 90087
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 3298
$ws.Range("F17").Value = "Trådticka"
$ws.Range("G17").Value = "Climacocystis borealis"
$ws.Range("H17").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q17").Value = 448761.1990147882
$ws.Range("R17").Value = 7087578.827763715

# Row 18
$ws.Range("A18").Value = 111936858
$ws.Range("B18").Value = 89845
$ws.Range("D18").Value = "VU"
$ws.Range("E18").Value = 1209
$ws.Range("F18").Value = "Rynkskinn"
$ws.Range("G18").Value = "Phlebia centrifuga"
$ws.Range("H18").Value = "P.Karst."
$ws.Range("Q18").Value = 448737.3665225056
$ws.Range("R18").Value = 7087496.445579056

# Row 19
$ws.Range("A19").Value = 111936867
$ws.Range("B19").Value = 89423
$ws.Range("E19").Value = 5432
$ws.Range("F19").Value = "Granticka"
$ws.Range("G19").Value = "Porodaedalea chrysoloma"
$ws.Range("H19").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("K19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("Q19").Value = 448791.554596175
$ws.Range("R19").Value = 7087386.366048628
$ws.Range("AC19").ClearContents()

# Row 20
$ws.Range("A20").Value = 111936796
$ws.Range("B20").Value = 56398
$ws.Range("E20").Value = 100109
$ws.Range("F20").Value = "Tretåig hackspett"
$ws.Range("G20").Value = "Picoides tridactylus"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("K20").Value = ""
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""
$ws.Range("Q20").Value = 448882.8980770012
$ws.Range("R20").Value = 7087229.443335658
$ws.Range("AC20").Value = "ringhack äldre"

# Row 21
$ws.Range("A21").Value = 111936795
$ws.Range("B21").Value = 56398
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = "Tretåig hackspett"
$ws.Range("G21").Value = "Picoides tridactylus"
$ws.Range("H21").Value = "(Linnaeus, 1758)"
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = ""
$ws.Range("N21").Value = ""
$ws.Range("Q21").Value = 448749.3706757246
$ws.Range("R21").Value = 7087421.839990681
$ws.Range("AC21").Value = "ringhack äldre"
